# Saldo.xlsx edit: account 004480134 (JOSE) moves from a small trailing
# balance (82.33) up near the top of the list with a much larger balance
# (41082.33), and the 004363260 (LARISSA) row is dropped entirely.
#
# Net row-count: 354 data rows -> 353 data rows (one row added, two removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Operate bottom-to-top so row numbers used below stay valid.

# 1) Remove the old "004480134 / JOSE / 82.33" row (currently row 180).
$ws.Rows.Item(180).Delete() | Out-Null

# 2) Remove the "004363260 / LARISSA / 4622.45" row (currently row 12).
$ws.Rows.Item(12).Delete() | Out-Null

# 3) Insert a new row above row 4 (currently "008035153 / CLAUDIO / 30051.49")
#    for the relocated "004480134 / JOSE" account with its new balance.
$ws.Rows.Item(4).Insert() | Out-Null

# Force the account number into the new row as literal text (matching the
# leading-zero text formatting used by every other row in column A) instead
# of letting it auto-convert to a number, then drop the temporary number
# format again so no stray style is left behind on the cell.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value2 = "004480134"
$ws.Range("A4").ClearFormats() | Out-Null

$ws.Range("B4").Value2 = "JOSE"
$ws.Range("C4").Value2 = 41082.33

Write-Host "Row 4:" $ws.Range("A4").Value2 $ws.Range("B4").Value2 $ws.Range("C4").Value2
Write-Host "Used rows:" $ws.UsedRange.Rows.Count
